# Add two new columns (I and J) with header row formatting matching the
# existing header style (bold, bordered, centered) used by column H ("IP").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header cell's formatting onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New header labels.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for rows 2-5.
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 7

$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 2
